$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.158
$ws.Range("E2").Value = 0.148
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 301
$ws.Range("L2").Value = 0.3266767961797265
$ws.Range("M2").Value = 78.122
$ws.Range("N2").Value = 0.07018417033510017
$ws.Range("O2").Value = 0.2595415282392027
$ws.Range("P2").Value = 78.09800000000001
$ws.Range("Q2").Value = 0.07016260893001529
$ws.Range("R2").Value = 0.2594617940199336
$ws.Range("S2").Value = 0.02399999999999913
$ws.Range("T2").Value = 0.0003072117969329911
$ws.Range("U2").Value = 1622.2
$ws.Range("V2").Value = 1.4573713053634
$ws.Range("W2").Value = 0.279407025459233
$ws.Range("X2").Value = 0.07401808066656144
$ws.Range("Y2").Value = 0.2053889447926716
$ws.Range("Z2").Value = 2.373518804739824
$ws.Range("AB2").Value = 0.07207246379987498
$ws.Range("AC2").Value = -0.07207246379987498
$ws.Range("AD2").Value = 565.7
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 565.7
$ws.Range("AG2").Value = -1056.5
$ws.Range("AH2").Value = 0.3369668811055516
$ws.Range("AI2").Value = 0.3048773915386688
$ws.Range("AJ2").Value = -18.66607773851593
$ws.Range("AK2").Value = -4.528504072010288

# Row 3
$ws.Range("B3").Value = "Standard Chartered Bank Ghana Limited (GHSE:SCB)"
$ws.Range("D3").Value = 0.162
$ws.Range("E3").Value = 0.212
$ws.Range("K3").Value = 73.8
$ws.Range("L3").Value = 0.4336075205640423
$ws.Range("M3").Value = 36.7
$ws.Range("N3").Value = 0.09721854304635762
$ws.Range("O3").Value = 0.497289972899729
$ws.Range("P3").Value = 36.7
$ws.Range("Q3").Value = 0.09721854304635762
$ws.Range("R3").Value = 0.497289972899729
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 557.5
$ws.Range("V3").Value = 1.47682119205298
$ws.Range("W3").Value = 0.3614103819784525
$ws.Range("X3").Value = 0.06378551655920309
$ws.Range("Y3").Value = 0.2976248654192494
$ws.Range("Z3").Value = -2.103831891223732
$ws.Range("AB3").Value = 0.06353193794611091
$ws.Range("AC3").Value = -0.06353193794611091
$ws.Range("AD3").Value = 15.9
$ws.Range("AF3").Value = 15.9
$ws.Range("AG3").Value = -541.6
$ws.Range("AH3").Value = 0.04041687849517031
$ws.Range("AI3").Value = 0.06382978723404256
$ws.Range("AJ3").Value = 3.300426569165143
$ws.Range("AK3").Value = 1.756160830090791

# Row 4
$ws.Range("B4").Value = "Ecobank Ghana Limited (GHSE:EGH)"
$ws.Range("D4").Value = 0.0975
$ws.Range("E4").Value = 0.0545
$ws.Range("K4").Value = 86.7
$ws.Range("L4").Value = 0.3244760479041917
$ws.Range("M4").Value = 16.7752
$ws.Range("N4").Value = 0.04205364753070946
$ws.Range("O4").Value = 0.1934855824682814
$ws.Range("P4").Value = 16.7752
$ws.Range("Q4").Value = 0.04205364753070946
$ws.Range("R4").Value = 0.1934855824682814
$ws.Range("U4").Value = 436.7
$ws.Range("V4").Value = 1.094760591626974
$ws.Range("W4").Value = 0.279407025459233
$ws.Range("X4").Value = 0.06468905854663469
$ws.Range("Y4").Value = 0.2147179669125983
$ws.Range("Z4").Value = -3.078341013824886
$ws.Range("AA4").Value = -0
$ws.Range("AB4").Value = 0.06474438143933785
$ws.Range("AC4").Value = -0.06474438143933785
$ws.Range("AD4").Value = 25.9
$ws.Range("AF4").Value = 25.9
$ws.Range("AG4").Value = -410.8
$ws.Range("AH4").Value = 0.06096986817325801
$ws.Range("AI4").Value = 0.06671818650180319
$ws.Range("AJ4").Value = 34.52100840336124
$ws.Range("AK4").Value = 8.470103092783505

# Row 5
$ws.Range("B5").Value = "GCB Bank Limited (GHSE:GCB)"
$ws.Range("D5").Value = 0.194
$ws.Range("E5").Value = 0.148
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 85
$ws.Range("L5").Value = 0.2812706816677697
$ws.Range("M5").Value = 9.24
$ws.Range("N5").Value = 0.05010845986984815
$ws.Range("O5").Value = 0.1087058823529412
$ws.Range("P5").Value = 9.24
$ws.Range("Q5").Value = 0.05010845986984815
$ws.Range("R5").Value = 0.1087058823529412
$ws.Range("U5").Value = 383.2
$ws.Range("V5").Value = 2.078091106290672
$ws.Range("W5").Value = 0.2890173410404624
$ws.Range("X5").Value = 0.07401808066656144
$ws.Range("Y5").Value = 0.2149992603739009
$ws.Range("Z5").Value = -21.28169014084517
$ws.Range("AA5").Value = -0
$ws.Range("AB5").Value = 0.07207246379987498
$ws.Range("AC5").Value = -0.07207246379987498
$ws.Range("AD5").Value = 55.4
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 55.4
$ws.Range("AG5").Value = -327.8
$ws.Range("AH5").Value = 0.2310258548790659
$ws.Range("AI5").Value = 0.1354854487649792
$ws.Range("AJ5").Value = 2.285913528591353
$ws.Range("AK5").Value = -12.75486381322958

# Row 6
$ws.Range("B6").Value = "Societe Generale Ghana Limited (GHSE:SOGEGH)"
$ws.Range("D6").Value = 0.158
$ws.Range("E6").Value = 0.237
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 24.3
$ws.Range("L6").Value = 0.2825581395348837
$ws.Range("M6").Value = 5.672800000000001
$ws.Range("N6").Value = 0.07272820512820513
$ws.Range("O6").Value = 0.2334485596707819
$ws.Range("P6").Value = 5.672800000000001
$ws.Range("Q6").Value = 0.07272820512820513
$ws.Range("R6").Value = 0.2334485596707819
$ws.Range("U6").Value = 156.5
$ws.Range("V6").Value = 2.006410256410256
$ws.Range("W6").Value = 0.1745689655172414
$ws.Range("X6").Value = 0.1024408489514024
$ws.Range("Y6").Value = 0.07212811656583899
$ws.Range("Z6").Value = 0.7543859649122809
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.08385478722217873
$ws.Range("AC6").Value = -0.08385478722217873
$ws.Range("AD6").Value = 79.40000000000001
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 79.40000000000001
$ws.Range("AG6").Value = -77.09999999999999
$ws.Range("AH6").Value = 0.5044472681067345
$ws.Range("AI6").Value = 0.3371549893842888
$ws.Range("AJ6").Value = -85.66666666666612
$ws.Range("AK6").Value = -0.9759493670886076

# Row 7
$ws.Range("B7").Value = "CalBank PLC (GHSE:CAL)"
$ws.Range("D7").Value = 0.09119999999999999
$ws.Range("E7").Value = 0.0272
$ws.Range("K7").Value = 31.2
$ws.Range("L7").Value = 0.325678496868476
$ws.Range("M7").Value = 9.734
$ws.Range("N7").Value = 0.1310094212651413
$ws.Range("O7").Value = 0.3119871794871795
$ws.Range("P7").Value = 9.710000000000001
$ws.Range("Q7").Value = 0.1306864064602961
$ws.Range("R7").Value = 0.3112179487179488
$ws.Range("S7").Value = 0.02399999999999913
$ws.Range("T7").Value = 0.002465584549003404
$ws.Range("U7").Value = 88.3
$ws.Range("V7").Value = 1.18842530282638
$ws.Range("W7").Value = 0.1894353369763206
$ws.Range("X7").Value = 0.2695644019956929
$ws.Range("Y7").Value = -0.08012906501937234
$ws.Range("Z7").Value = 0.2100416575312432
$ws.Range("AB7").Value = 0.09151129941650554
$ws.Range("AC7").Value = -0.09151129941650554
$ws.Range("AD7").Value = 389.1
$ws.Range("AF7").Value = 389.1
$ws.Range("AG7").Value = 300.8
$ws.Range("AH7").Value = 0.839663357790246
$ws.Range("AI7").Value = 0.6781108400139423
$ws.Range("AJ7").Value = 0.8019194881364969
$ws.Range("AK7").Value = 0.61956745623069

# Remove cells no longer present in the updated dataset
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()
